$d = $word.ActiveDocument

# --- Table 2 ("Command Codes"): append row 30 | Set Stepper 1 Command ---
$tCommandCodes = $d.Tables.Item(2)
$row1 = $tCommandCodes.Rows.Add()
$row1.Cells.Item(1).Range.Text = "30"
$row1.Cells.Item(2).Range.Text = "Set Stepper 1 Command"

# --- Table 3 ("Mapping Codes"): append row Set Stepper X Command | Speed (rpm) ---
$tMappingCodes = $d.Tables.Item(3)
$row2 = $tMappingCodes.Rows.Add()
$row2.Cells.Item(1).Range.Text = "Set Stepper X Command"
$row2.Cells.Item(2).Range.Text = "Speed (rpm)"

# --- Table 4 ("State Codes"): append row Set Stepper X Command | Steps ---
$tStateCodes = $d.Tables.Item(4)
$row3 = $tStateCodes.Rows.Add()
$row3.Cells.Item(1).Range.Text = "Set Stepper X Command"
$row3.Cells.Item(2).Range.Text = "Steps"
